$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wong3")
$ws.Name = "Euclid"

$ws.Range("C2").Value = 403
$ws.Range("D2").Value = 89.95535714285714
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0.2232142857142857
$ws.Range("C4").Value = 347
$ws.Range("D4").Value = 77.45535714285714
$ws.Range("C5").Value = 229
$ws.Range("D5").Value = 51.11607142857143
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 0.6696428571428571
$ws.Range("C7").Value = 327
$ws.Range("D7").Value = 72.99107142857143
$ws.Range("C8").Value = 395
$ws.Range("D8").Value = 88.16964285714286
$ws.Range("C9").Value = 224
$ws.Range("D9").Value = 50
$ws.Range("C10").Value = 213
$ws.Range("D10").Value = 47.54464285714285
$ws.Range("C11").Value = 225
$ws.Range("D11").Value = 50.22321428571429
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 0.4464285714285714
$ws.Range("C13").Value = 326
$ws.Range("D13").Value = 72.76785714285714
$ws.Range("C14").Value = 405
$ws.Range("D14").Value = 90.40178571428571
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 0.4464285714285714
$ws.Range("C16").Value = 350
$ws.Range("D16").Value = 78.125
$ws.Range("C17").Value = 201
$ws.Range("D17").Value = 44.86607142857143
$ws.Range("C18").Value = 399
$ws.Range("D18").Value = 89.0625
$ws.Range("C19").Value = 420
$ws.Range("D19").Value = 93.75
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 0.2232142857142857
$ws.Range("C21").Value = 327
$ws.Range("D21").Value = 72.99107142857143
$ws.Range("C22").Value = 149
$ws.Range("D22").Value = 33.25892857142857
$ws.Range("C23").Value = 231
$ws.Range("D23").Value = 51.5625
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 0.2232142857142857
$ws.Range("C25").Value = 249
$ws.Range("D25").Value = 55.58035714285714
$ws.Range("C26").Value = 420
$ws.Range("D26").Value = 93.75
$ws.Range("C27").Value = 233
$ws.Range("D27").Value = 52.00892857142857
$ws.Range("C28").Value = 403
$ws.Range("D28").Value = 89.95535714285714
$ws.Range("C29").Value = 3
$ws.Range("D29").Value = 0.6696428571428571
$ws.Range("C30").Value = 395
$ws.Range("D30").Value = 88.16964285714286
$ws.Range("C31").Value = 420
$ws.Range("D31").Value = 93.75
$ws.Range("C32").Value = 391
$ws.Range("D32").Value = 87.27678571428571
$ws.Range("C33").Value = 213
$ws.Range("D33").Value = 47.54464285714285
$ws.Range("C35").Value = 5
$ws.Range("D35").Value = 1.116071428571429
$ws.Range("C36").Value = 347
$ws.Range("D36").Value = 77.45535714285714
$ws.Range("C37").Value = 148
$ws.Range("D37").Value = 33.03571428571428
$ws.Range("C38").Value = 399
$ws.Range("D38").Value = 89.0625
$ws.Range("C39").Value = 229
$ws.Range("D39").Value = 51.11607142857143
$ws.Range("C40").Value = 397
$ws.Range("D40").Value = 88.61607142857143
$ws.Range("C41").Value = 321
$ws.Range("D41").Value = 71.65178571428571
$ws.Range("C42").Value = 399
$ws.Range("D42").Value = 89.0625
$ws.Range("C43").Value = 148
$ws.Range("D43").Value = 33.03571428571428
$ws.Range("C44").Value = 346
$ws.Range("D44").Value = 77.23214285714286
$ws.Range("C45").Value = 1
$ws.Range("D45").Value = 0.2232142857142857
$ws.Range("C46").Value = 420
$ws.Range("D46").Value = 93.75
$ws.Range("C47").Value = 143
$ws.Range("D47").Value = 31.91964285714285
$ws.Range("C48").Value = 346
$ws.Range("D48").Value = 77.23214285714286
$ws.Range("C49").Value = 3
$ws.Range("D49").Value = 0.6696428571428571
$ws.Range("C50").Value = 348
$ws.Range("D50").Value = 77.67857142857143
$ws.Range("C51").Value = 318
$ws.Range("D51").Value = 70.98214285714286
$ws.Range("C52").Value = 347
$ws.Range("D52").Value = 77.45535714285714
$ws.Range("C53").Value = 420
$ws.Range("D53").Value = 93.75
$ws.Range("C54").Value = 401
$ws.Range("D54").Value = 89.50892857142857
$ws.Range("C55").Value = 319
$ws.Range("D55").Value = 71.20535714285714
$ws.Range("C56").Value = 349
$ws.Range("D56").Value = 77.90178571428571
$ws.Range("C57").Value = 2
$ws.Range("D57").Value = 0.4464285714285714
$ws.Range("C58").Value = 395
$ws.Range("D58").Value = 88.16964285714286
$ws.Range("C59").Value = 151
$ws.Range("D59").Value = 33.70535714285715
$ws.Range("C60").Value = 60
$ws.Range("D60").Value = 13.39285714285714
$ws.Range("C61").Value = 304
$ws.Range("D61").Value = 67.85714285714286
$ws.Range("C63").Value = 4
$ws.Range("D63").Value = 0.8928571428571428
$ws.Range("C65").Value = 3
$ws.Range("D65").Value = 0.6696428571428571
$ws.Range("C66").Value = 409
$ws.Range("D66").Value = 91.29464285714286
$ws.Range("C67").Value = 3
$ws.Range("D67").Value = 0.6696428571428571
$ws.Range("C68").Value = 3
$ws.Range("D68").Value = 0.6696428571428571
$ws.Range("C69").Value = 249
$ws.Range("D69").Value = 55.58035714285714
$ws.Range("C70").Value = 181
$ws.Range("D70").Value = 40.40178571428572
$ws.Range("C71").Value = 71
$ws.Range("D71").Value = 15.84821428571428
$ws.Range("C72").Value = 409
$ws.Range("D72").Value = 91.29464285714286
$ws.Range("C73").Value = 1
$ws.Range("D73").Value = 0.2232142857142857
$ws.Range("C74").Value = 420
$ws.Range("D74").Value = 93.75
$ws.Range("C75").Value = 322
$ws.Range("D75").Value = 71.875
$ws.Range("C76").Value = 401
$ws.Range("D76").Value = 89.50892857142857
$ws.Range("C77").Value = 196
$ws.Range("D77").Value = 43.75
$ws.Range("C78").Value = 407
$ws.Range("D78").Value = 90.84821428571429
$ws.Range("C79").Value = 2
$ws.Range("D79").Value = 0.4464285714285714
$ws.Range("C81").Value = 259
$ws.Range("D81").Value = 57.8125
$ws.Range("C82").Value = 157
$ws.Range("D82").Value = 35.04464285714285
$ws.Range("C83").Value = 322
$ws.Range("D83").Value = 71.875
